# DEV: Change the project names of data/delete category.
# Renames the "data/delete" group of services (rows 17-20) from the old
# "selling" sale-delete endpoints to the new "orderizer" orders-delete
# endpoints, and un-bolds rows 14-16 (data/save group) to match the rest
# of the non-header rows in that block. Also moves the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C (service names) filled top-to-bottom first ---
$ws.Range("C17").Value = "data-delete-orders-manager"
$ws.Range("C18").Value = "data-delete-order"
$ws.Range("C19").Value = "data-delete-orders"
$ws.Range("C20").Value = "data-delete-search-orders"

# --- Column E (base paths): rows 18-20 filled first, row 17 (manager) last ---
$ws.Range("E18").Value = " /orderizer/data/delete/order/v1/"
$ws.Range("E19").Value = " /orderizer/data/delete/orders/v1/"
$ws.Range("E20").Value = " /orderizer/data/delete/search/orders/v1/"
$ws.Range("E17").Value = "/orderizer/data/delete/orders/manager/v1/"

# --- Rows 14-16: drop the bold emphasis so they match the plain fill used
#     elsewhere in the sheet (style index 11 -> 5) ---
$ws.Range("A14:E16").Font.Bold = $false

# --- Move the active cell/selection like the saved workbook view ---
$ws.Range("F12").Select()
